$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the "Price" column keeps its original text formatting (values like
# "30.575.23" or "1.000" must not be auto-converted to numbers by Excel).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.575.23"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.872.41"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.21"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2906"
$ws.Range("E8").Value = "  +1.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06475"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.13"
$ws.Range("E10").Value = "  +4.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07711"
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7419"
$ws.Range("E12").Value = "  +4.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.53"
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.868.85"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.148"
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "273.43"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.559.18"
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.34"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9998"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007493"
$ws.Range("E20").Value = "  -0.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.116.54"
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.254"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.188"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.232"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.33"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.79"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.916"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09980"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.505"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.293"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.106"
$ws.Range("E33").Value = "  +1.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04784"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.118"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6963"
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9999"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.716"
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01848"
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.751"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.189"
$ws.Range("E41").Value = "  -1.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.43"
$ws.Range("E42").Value = "  +4.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.967"
$ws.Range("E43").Value = "  +2.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4178"
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.8330"
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.65"
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.301"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.37"
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "926.08"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.959"
$ws.Range("E51").Value = "  -1.76%  "
